$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 412.92
$ws.Range("I9").Value = 236.52632
$ws.Range("K9").Value = 236.52632
$ws.Range("M9").Value = -67.52632

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 3849.5557
$ws.Range("I38").Value = 149.33333
$ws.Range("J38").Value = 11250
$ws.Range("K38").Value = 447.99999
$ws.Range("L38").Value = 33750
$ws.Range("M38").Value = -75.99998999999997
$ws.Range("N38").Value = -34494

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 46883110
$ws.Range("I62").Value = 19238638
$ws.Range("K62").Value = 19238638
$ws.Range("M62").Value = -19238014

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 6852.2
$ws.Range("I64").Value = 5151.4287
$ws.Range("J64").Value = 7768
$ws.Range("K64").Value = 5151.4287
$ws.Range("L64").Value = 7768
$ws.Range("M64").Value = -4903.4287
$ws.Range("N64").Value = -8264

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 46883110
$ws.Range("I65").Value = 19238638
$ws.Range("K65").Value = 96193190
$ws.Range("M65").Value = -96190070

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 6852.2
$ws.Range("I67").Value = 5151.4287
$ws.Range("J67").Value = 7768
$ws.Range("K67").Value = 5151.4287
$ws.Range("L67").Value = 7768
$ws.Range("M67").Value = -4293.4287
$ws.Range("N67").Value = -9484

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3581.9524
$ws.Range("I132").Value = 2866.8333
$ws.Range("K132").Value = 8600.499899999999
$ws.Range("M132").Value = -6070.499899999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5030.19
$ws.Range("J138").Value = 5679.577
$ws.Range("L138").Value = 17038.731
$ws.Range("N138").Value = -27318.731

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6363.4194
$ws.Range("I2").Value = 6423.5713
$ws.Range("J2").Value = 6237.1
$ws.Range("K2").Value = 6423.5713
$ws.Range("L2").Value = 6237.1
$ws.Range("M2").Value = -6310.5713
$ws.Range("N2").Value = -6463.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3425.182
$ws.Range("I61").Value = 2119.2856
$ws.Range("K61").Value = 2119.2856
$ws.Range("M61").Value = -1907.2856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4261.4385
$ws.Range("I74").Value = 4258.75
$ws.Range("J74").Value = 4298
$ws.Range("K74").Value = 4258.75
$ws.Range("L74").Value = 4298
$ws.Range("M74").Value = -3384.75
$ws.Range("N74").Value = -6046

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4261.4385
$ws.Range("I77").Value = 4258.75
$ws.Range("J77").Value = 4298
$ws.Range("K77").Value = 21293.75
$ws.Range("L77").Value = 21490
$ws.Range("M77").Value = -16925.75
$ws.Range("N77").Value = -30226

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1707.1
$ws.Range("I88").Value = 1675.5
$ws.Range("J88").Value = 1720.6428
$ws.Range("K88").Value = 1675.5
$ws.Range("L88").Value = 1720.6428
$ws.Range("M88").Value = -1269.5
$ws.Range("N88").Value = -2532.6428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1707.1
$ws.Range("I91").Value = 1675.5
$ws.Range("J91").Value = 1720.6428
$ws.Range("K91").Value = 1675.5
$ws.Range("L91").Value = 1720.6428
$ws.Range("M91").Value = -271.5
$ws.Range("N91").Value = -4528.6428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2998.7222
$ws.Range("I102").Value = 2461.0833
$ws.Range("J102").Value = 4074
$ws.Range("K102").Value = 2461.0833
$ws.Range("L102").Value = 4074
$ws.Range("M102").Value = -839.0832999999998
$ws.Range("N102").Value = -7318

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 6363.4194
$ws.Range("I116").Value = 6423.5713
$ws.Range("J116").Value = 6237.1
$ws.Range("K116").Value = 6423.5713
$ws.Range("L116").Value = 6237.1
$ws.Range("M116").Value = -4129.5713
$ws.Range("N116").Value = -10825.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 7414.143
$ws.Range("I132").Value = 6480
$ws.Range("J132").Value = 9749.5
$ws.Range("K132").Value = 19440
$ws.Range("L132").Value = 29248.5
$ws.Range("M132").Value = -16910
$ws.Range("N132").Value = -34308.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 79000
$ws.Range("J133").Value = 79000
$ws.Range("L133").Value = 79000
$ws.Range("N133").Value = -84060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3425.182
$ws.Range("I136").Value = 2119.2856
$ws.Range("K136").Value = 6357.8568
$ws.Range("M136").Value = -3807.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6363.4194
$ws.Range("I3").Value = 6423.5713
$ws.Range("J3").Value = 6237.1
$ws.Range("K3").Value = 6423.5713
$ws.Range("L3").Value = 6237.1
$ws.Range("M3").Value = -6309.5713
$ws.Range("N3").Value = -6465.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2044.7778
$ws.Range("I86").Value = 1550.5
$ws.Range("K86").Value = 1550.5
$ws.Range("M86").Value = -427.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2044.7778
$ws.Range("I89").Value = 1550.5
$ws.Range("K89").Value = 7752.5
$ws.Range("M89").Value = -2136.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2777.8572
$ws.Range("I134").Value = 3472.5625
$ws.Range("K134").Value = 10417.6875
$ws.Range("M134").Value = -7882.6875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3170.0645
$ws.Range("I132").Value = 2492.5925
$ws.Range("K132").Value = 7477.7775
$ws.Range("M132").Value = -4947.7775

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 600.5714
$ws.Range("I46").Value = 495
$ws.Range("J46").Value = 679.75
$ws.Range("K46").Value = 1485
$ws.Range("L46").Value = 2039.25
$ws.Range("M46").Value = -1394
$ws.Range("N46").Value = -2221.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3847.4
$ws.Range("I102").Value = 5149.5
$ws.Range("K102").Value = 5149.5
$ws.Range("M102").Value = -3527.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 53383.95
$ws.Range("I132").Value = 70241.87
$ws.Range("K132").Value = 210725.61
$ws.Range("M132").Value = -208195.61

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 48735.668
$ws.Range("J134").Value = 48735.668
$ws.Range("L134").Value = 146207.004
$ws.Range("N134").Value = -151277.004

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2468.7
$ws.Range("I46").Value = 1520.7778
$ws.Range("J46").Value = 11000
$ws.Range("K46").Value = 1520.7778
$ws.Range("L46").Value = 11000
$ws.Range("M46").Value = -1332.7778
$ws.Range("N46").Value = -11376

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 12822216
$ws.Range("I132").Value = 18519932
$ws.Range("J132").Value = 2351.875
$ws.Range("K132").Value = 55559796
$ws.Range("L132").Value = 7055.625
$ws.Range("M132").Value = -55557266

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3737.3296
$ws.Range("J136").Value = 4594.4
$ws.Range("L136").Value = 13783.2
$ws.Range("N136").Value = -18883.2
